$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0) Remove the existing table definition so we can freely restructure the
#    range beneath it; we will recreate "Tableau1" once the grid is final.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()

# ---------------------------------------------------------------------------
# 1) Insert a new column at D (shift existing D:H right to E:I) to make room
#    for the new "PÉRIODE" column.
# ---------------------------------------------------------------------------
$ws.Range("D1:D6").Insert(-4161)   # xlShiftToRight

# ---------------------------------------------------------------------------
# 2) Rename "PRIORITÉ" (col A) to "CODE" and replace its numeric priority
#    values with the new text codes.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "CODE"
$ws.Range("A2").Value = "UICN"
$ws.Range("A3").Value = "SLL"
$ws.Range("A4").Value = "SLL"
$ws.Range("A5").Value = "OPP"
$ws.Range("A6").Clear()

# ---------------------------------------------------------------------------
# 3) Fill in the new "PÉRIODE" column (values first, number format after, so
#    the numeric 2016 stays a real number instead of being coerced to text).
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "PÉRIODE"
$ws.Range("D2").Value = "2014-2016, 2018, 2020-2023"
$ws.Range("D3").Value = 2016
$ws.Range("D4").Value = "2013-2014"
$ws.Range("D5").Value = "2011-2023"
$ws.Range("D6").Clear()

$ws.Range("D1:D5").NumberFormat = "@"
$ws.Range("D1:D5").WrapText = $true
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# ---------------------------------------------------------------------------
# 4) Recreate the table over the new extent with the columns in their
#    (now-correct) physical order.
# ---------------------------------------------------------------------------
$newTable = $ws.ListObjects.Add(1, $ws.Range("A1:I6"), 0, 1)
$newTable.Name = "Tableau1"
$newTable.TableStyle = "TableStyleMedium13"

$newTable.ListColumns.Item(1).Name = "CODE"
$newTable.ListColumns.Item(2).Name = "PROTOCOLE"
$newTable.ListColumns.Item(3).Name = "TYPE DE PROTOCOLE"
$newTable.ListColumns.Item(4).Name = "PÉRIODE"
$newTable.ListColumns.Item(5).Name = "N TOTAL"
$newTable.ListColumns.Item(6).Name = "N PRÉSENCE"
$newTable.ListColumns.Item(7).Name = "MOTS CLÉS (| = ""OU"" ;  & = ""ET"" ;  ! = ""NON"" ; X = chiffre)"
$newTable.ListColumns.Item(8).Name = "COMMENTAIRE"
$newTable.ListColumns.Item(9).Name = "DESCRIPTION (OPTIONNEL)"

# ---------------------------------------------------------------------------
# 5) Row heights for the two wrapped comment rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45

$ws.Range("F21").Select()
